$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 329.5
$ws.Range("I2").Value = 329.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 329.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -216.5
$ws.Range("N2").Value = ""
$ws.Range("H19").Value = 1446.4546
$ws.Range("I19").Value = 1576
$ws.Range("J19").Value = 1291
$ws.Range("K19").Value = 1576
$ws.Range("L19").Value = 1291
$ws.Range("M19").Value = -1401
$ws.Range("N19").Value = -1641
$ws.Range("H31").Value = 229
$ws.Range("I31").Value = 229
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 687
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -457
$ws.Range("N31").Value = ""
$ws.Range("H32").Value = 1818
$ws.Range("J32").Value = 1999.3334
$ws.Range("L32").Value = 1999.3334
$ws.Range("N32").Value = -2651.3334
$ws.Range("H64").Value = 6844
$ws.Range("I64").Value = 6844
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 6844
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -6596
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 6844
$ws.Range("I67").Value = 6844
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 6844
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -5986
$ws.Range("N67").Value = ""
$ws.Range("H113").Value = 12000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -18508
$ws.Range("H137").Value = 23811380
$ws.Range("I137").Value = 83334710
$ws.Range("K137").Value = 250004130
$ws.Range("M137").Value = -250001580
$ws.Range("N137").Value = ""
$ws.Range("H138").Value = 3571.3333
$ws.Range("I138").Value = 3904.8572
$ws.Range("J138").Value = 3325.5789
$ws.Range("K138").Value = 11714.5716
$ws.Range("L138").Value = 9976.736699999999
$ws.Range("M138").Value = -6574.571599999999
$ws.Range("N138").Value = -20256.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1657.5
$ws.Range("I29").Value = 1657.5
$ws.Range("K29").Value = 1657.5
$ws.Range("M29").Value = -1349.5
$ws.Range("H32").Value = 8772.521000000001
$ws.Range("I32").Value = 5654.1113
$ws.Range("K32").Value = 5654.1113
$ws.Range("M32").Value = -5367.1113
$ws.Range("N32").Value = ""
$ws.Range("H74").Value = 2129.5334
$ws.Range("I74").Value = 1924.5714
$ws.Range("K74").Value = 1924.5714
$ws.Range("M74").Value = -1050.5714
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 2129.5334
$ws.Range("I77").Value = 1924.5714
$ws.Range("K77").Value = 9622.857
$ws.Range("M77").Value = -5254.857
$ws.Range("N77").Value = ""
$ws.Range("H132").Value = 3276.4614
$ws.Range("I132").Value = 3276.4614
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9829.3842
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7299.3842
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 812.53845
$ws.Range("I94").Value = 757.4
$ws.Range("J94").Value = 996.3333
$ws.Range("K94").Value = 757.4
$ws.Range("L94").Value = 996.3333
$ws.Range("M94").Value = -306.4
$ws.Range("N94").Value = -1898.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 373.6
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
$ws.Range("H25").Value = 3913.625
$ws.Range("I25").Value = 3913.625
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 3913.625
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -3739.625
$ws.Range("N25").Value = ""
$ws.Range("H86").Value = 15799.4
$ws.Range("I86").Value = 14000.286
$ws.Range("K86").Value = 14000.286
$ws.Range("M86").Value = -12877.286
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 15799.4
$ws.Range("I89").Value = 14000.286
$ws.Range("K89").Value = 70001.42999999999
$ws.Range("M89").Value = -64385.42999999999
$ws.Range("N89").Value = ""
$ws.Range("H134").Value = 2118.4285
$ws.Range("I134").Value = 2439.6
$ws.Range("J134").Value = 1315.5
$ws.Range("K134").Value = 7318.799999999999
$ws.Range("L134").Value = 3946.5
$ws.Range("M134").Value = -4783.799999999999
$ws.Range("N134").Value = -9016.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9803.866
$ws.Range("I56").Value = 9803.866
$ws.Range("K56").Value = 9803.866
$ws.Range("M56").Value = -9273.866
$ws.Range("H122").Value = 865.2857
$ws.Range("I122").Value = 864.25
$ws.Range("J122").Value = 866.6667
$ws.Range("K122").Value = 7778.25
$ws.Range("L122").Value = 7800.0003
$ws.Range("M122").Value = -5328.25
$ws.Range("N122").Value = -12700.0003
$ws.Range("H140").Value = 501467.8
$ws.Range("I140").Value = 501467.8
$ws.Range("K140").Value = 1504403.4
$ws.Range("M140").Value = -1499223.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 4900
$ws.Range("I24").Value = 4900
$ws.Range("K24").Value = 4900
$ws.Range("M24").Value = -4727
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744
$ws.Range("H126").Value = 6498.9375
$ws.Range("I126").Value = 6588.4546
$ws.Range("K126").Value = 19765.3638
$ws.Range("M126").Value = -17295.3638
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 18522648
$ws.Range("J132").Value = 55559424
$ws.Range("L132").Value = 166678272
$ws.Range("N132").Value = -166683332
$ws.Range("H136").Value = 69990.44500000001
$ws.Range("J136").Value = 69990.44500000001
$ws.Range("L136").Value = 209971.335
$ws.Range("N136").Value = -215071.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1500
$ws.Range("K7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("H22").Value = 45456350
$ws.Range("I22").Value = 995.6
$ws.Range("J22").Value = 83335816
$ws.Range("K22").Value = 995.6
$ws.Range("L22").Value = 83335816
$ws.Range("M22").Value = -700.6
$ws.Range("N22").Value = -83336406
$ws.Range("H27").Value = 45456350
$ws.Range("I27").Value = 995.6
$ws.Range("J27").Value = 83335816
$ws.Range("K27").Value = 995.6
$ws.Range("L27").Value = 83335816
$ws.Range("M27").Value = -888.6
$ws.Range("N27").Value = -83336030
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H122").Value = 3501.7144
$ws.Range("I122").Value = 3252
$ws.Range("K122").Value = 9756
$ws.Range("M122").Value = -7306
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H130").Value = 96623
$ws.Range("J130").Value = 96623
$ws.Range("L130").Value = 96623
$ws.Range("N130").Value = -106663
$ws.Range("H132").Value = 3706.4666
$ws.Range("I132").Value = 3949.5
$ws.Range("J132").Value = 3618.0908
$ws.Range("K132").Value = 11848.5
$ws.Range("L132").Value = 10854.2724
$ws.Range("M132").Value = -9318.5
$ws.Range("N132").Value = -15914.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2624.75
$ws.Range("J5").Value = 2624.75
$ws.Range("L5").Value = 2624.75
$ws.Range("N5").Value = -2848.75
$ws.Range("H100").Value = 1811.1111
$ws.Range("I100").Value = 860.2
$ws.Range("J100").Value = 2999.75
$ws.Range("K100").Value = 1720.4
$ws.Range("L100").Value = 5999.5
$ws.Range("M100").Value = -1179.4
$ws.Range("N100").Value = -7081.5
$ws.Range("H125").Value = 78198
$ws.Range("J125").Value = 78198
$ws.Range("L125").Value = 78198
$ws.Range("N125").Value = -88038
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
